# Weekly update: insert this week's two new price rows (Zafiro rojo / Zafiro
# verde, "Terminal Hortofrutícola Agro Chillán") at the top of the data block
# (row 62), pushing all previously-recorded rows down by two. The sheet's
# used range grows from A1:R153 to A1:R155 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 62 - everything that
# used to live on rows 62..153 shifts down to rows 64..155.
$ws.Rows.Item(62).Insert()
$ws.Rows.Item(62).Insert()

# New row 62: Zafiro rojo, Primera
$ws.Cells.Item(62, 1).Value  = 7
$ws.Cells.Item(62, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(62, 3).Value  = "Ñuble"
$ws.Cells.Item(62, 4).Value  = 44483
$ws.Cells.Item(62, 5).Value  = 16
$ws.Cells.Item(62, 6).Value  = 100112002
$ws.Cells.Item(62, 7).Value  = "Pimiento"
$ws.Cells.Item(62, 8).Value  = "Zafiro rojo"
$ws.Cells.Item(62, 9).Value  = "Primera"
$ws.Cells.Item(62, 10).Value = 160
$ws.Cells.Item(62, 11).Value = 43000
$ws.Cells.Item(62, 12).Value = 44000
$ws.Cells.Item(62, 13).Value = 43500
$ws.Cells.Item(62, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 2900
$ws.Cells.Item(62, 17).Value = 15
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# New row 63: Zafiro verde, Primera
$ws.Cells.Item(63, 1).Value  = 7
$ws.Cells.Item(63, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value  = "Ñuble"
$ws.Cells.Item(63, 4).Value  = 44483
$ws.Cells.Item(63, 5).Value  = 16
$ws.Cells.Item(63, 6).Value  = 100112002
$ws.Cells.Item(63, 7).Value  = "Pimiento"
$ws.Cells.Item(63, 8).Value  = "Zafiro verde"
$ws.Cells.Item(63, 9).Value  = "Primera"
$ws.Cells.Item(63, 10).Value = 160
$ws.Cells.Item(63, 11).Value = 41000
$ws.Cells.Item(63, 12).Value = 42000
$ws.Cells.Item(63, 13).Value = 41500
$ws.Cells.Item(63, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(63, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(63, 16).Value = 2767
$ws.Cells.Item(63, 17).Value = 15
$ws.Cells.Item(63, 18).Value = "Hortaliza"
